$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-10-02"

# Update the October label in column A (row 11) to reflect the new date
$ws.Range("A11").Value = "October (through 10-02)"

# Fill in the October data (row 11) with the updated values
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 6
$ws.Range("D11").Value = 8
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 9
$ws.Range("H11").Value = 17
$ws.Range("I11").Value = 9

# Update the Total row (row 12) with the new totals
$ws.Range("B12").Value = 229
$ws.Range("C12").Value = 435
$ws.Range("D12").Value = 635
$ws.Range("E12").Value = 553
$ws.Range("F12").Value = 423
$ws.Range("G12").Value = 910
$ws.Range("H12").Value = 1264
$ws.Range("I12").Value = 1291
